$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (the guide-RNA record that is no longer
# reported upstream); all subsequent rows shift up by one.
$ws.Rows.Item(1).Delete()
